$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 21376
$ws.Range("I16").Value = 19999
$ws.Range("J16").Value = 21572.715
$ws.Range("K16").Value = 19999
$ws.Range("L16").Value = 21572.715
$ws.Range("M16").Value = -19769
$ws.Range("N16").Value = -22032.715

$ws.Range("H80").Value = 1515972
$ws.Range("I80").Value = 2525833
$ws.Range("J80").Value = 1180.3334
$ws.Range("K80").Value = 7577499
$ws.Range("L80").Value = 3541.0002
$ws.Range("M80").Value = -7576501
$ws.Range("N80").Value = -5537.0002

$ws.Range("H83").Value = 1515972
$ws.Range("I83").Value = 2525833
$ws.Range("J83").Value = 1180.3334
$ws.Range("K83").Value = 22732497
$ws.Range("L83").Value = 10623.0006
$ws.Range("M83").Value = -22727505
$ws.Range("N83").Value = -20607.0006

$ws.Range("H116").Value = 30761376
$ws.Range("I116").Value = 27895634
$ws.Range("J116").Value = 33340546
$ws.Range("K116").Value = 27895634
$ws.Range("L116").Value = 33340546
$ws.Range("M116").Value = -27892192
$ws.Range("N116").Value = -33347430

$ws.Range("H127").Value = 1964.8572
$ws.Range("I127").Value = 1790.8
$ws.Range("K127").Value = 5372.4
$ws.Range("M127").Value = -412.3999999999996

$ws.Range("H137").Value = 17243862
$ws.Range("I137").Value = 71430780
$ws.Range("J137").Value = 2567
$ws.Range("K137").Value = 214292340
$ws.Range("L137").Value = 7701
$ws.Range("M137").Value = -214289790
$ws.Range("N137").Value = -12801

$ws.Range("H138").Value = 2971.842
$ws.Range("I138").Value = 2401.2144
$ws.Range("J138").Value = 3304.7083
$ws.Range("K138").Value = 7203.6432
$ws.Range("L138").Value = 9914.124899999999
$ws.Range("M138").Value = -2063.6432
$ws.Range("N138").Value = -20194.1249

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2469.1667
$ws.Range("I45").Value = 2092
$ws.Range("K45").Value = 2092
$ws.Range("M45").Value = -1715

$ws.Range("H46").Value = 8508.5
$ws.Range("J46").Value = 8508.5
$ws.Range("L46").Value = 8508.5
$ws.Range("N46").Value = -9146.5

$ws.Range("H122").Value = 2874.4285
$ws.Range("I122").Value = 1811.6428
$ws.Range("K122").Value = 5434.928400000001
$ws.Range("M122").Value = -2984.928400000001

$ws.Range("H125").Value = 97799.8
$ws.Range("J125").Value = 99750
$ws.Range("L125").Value = 99750
$ws.Range("N125").Value = -109590

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 28158.37
$ws.Range("I99").Value = 1499.0834
$ws.Range("J99").Value = 73860
$ws.Range("K99").Value = 1499.0834
$ws.Range("L99").Value = 73860
$ws.Range("M99").Value = -1.083399999999983
$ws.Range("N99").Value = -76856

$ws.Range("H105").Value = 2114.5334
$ws.Range("I105").Value = 2076.5
$ws.Range("K105").Value = 2076.5
$ws.Range("M105").Value = -329.5

$ws.Range("H107").Value = 11119692
$ws.Range("I107").Value = 4727.543
$ws.Range("J107").Value = 50022068
$ws.Range("K107").Value = 4727.543
$ws.Range("L107").Value = 50022068
$ws.Range("M107").Value = -2807.543
$ws.Range("N107").Value = -50025908

$ws.Range("H134").Value = 2103.0667
$ws.Range("I134").Value = 1903.2858
$ws.Range("K134").Value = 5709.857400000001
$ws.Range("M134").Value = -3174.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 300
$ws.Range("K12").Value = 300
$ws.Range("M12").Value = -130

$ws.Range("H19").Value = 486.1111
$ws.Range("I19").Value = 515.625
$ws.Range("K19").Value = 515.625
$ws.Range("M19").Value = -345.625

$ws.Range("H24").Value = 486.1111
$ws.Range("I24").Value = 515.625
$ws.Range("K24").Value = 515.625
$ws.Range("M24").Value = -345.625

$ws.Range("H31").Value = 3065.5833
$ws.Range("I31").Value = 2454.7778
$ws.Range("J31").Value = 3152.8413
$ws.Range("K31").Value = 2454.7778
$ws.Range("L31").Value = 3152.8413
$ws.Range("M31").Value = -2159.7778
$ws.Range("N31").Value = -3742.8413

$ws.Range("H34").Value = 3065.5833
$ws.Range("I34").Value = 2454.7778
$ws.Range("J34").Value = 3152.8413
$ws.Range("K34").Value = 2454.7778
$ws.Range("L34").Value = 3152.8413
$ws.Range("M34").Value = -2252.7778
$ws.Range("N34").Value = -3556.8413

$ws.Range("H99").Value = 2778.3845
$ws.Range("I99").Value = 2328.1428
$ws.Range("J99").Value = 3303.6667
$ws.Range("K99").Value = 2328.1428
$ws.Range("L99").Value = 3303.6667
$ws.Range("M99").Value = -830.1428000000001
$ws.Range("N99").Value = -6299.6667

$ws.Range("H126").Value = 2778.3845
$ws.Range("I126").Value = 2328.1428
$ws.Range("J126").Value = 3303.6667
$ws.Range("K126").Value = 6984.428400000001
$ws.Range("L126").Value = 9911.000100000001
$ws.Range("M126").Value = -4514.428400000001
$ws.Range("N126").Value = -14851.0001

$ws.Range("H132").Value = 2604.24
$ws.Range("J132").Value = 6160.6
$ws.Range("L132").Value = 18481.8
$ws.Range("N132").Value = -23541.8

$ws.Range("H134").Value = 2471.8708
$ws.Range("I134").Value = 1738.6842
$ws.Range("K134").Value = 5216.0526
$ws.Range("M134").Value = -2681.0526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 673.875
$ws.Range("I5").Value = 324.33334
$ws.Range("K5").Value = 973.0000200000001
$ws.Range("M5").Value = -861.0000200000001

$ws.Range("H23").Value = 1310.7273
$ws.Range("I23").Value = 775
$ws.Range("K23").Value = 2325
$ws.Range("M23").Value = -2090

$ws.Range("H33").Value = 189.375
$ws.Range("I33").Value = 155
$ws.Range("K33").Value = 930
$ws.Range("M33").Value = -647

$ws.Range("H68").Value = 1817.3529
$ws.Range("J68").Value = 1938.9286
$ws.Range("L68").Value = 5816.7858
$ws.Range("N68").Value = -7438.7858

$ws.Range("H71").Value = 1817.3529
$ws.Range("J71").Value = 1938.9286
$ws.Range("L71").Value = 17450.3574
$ws.Range("N71").Value = -25562.3574

$ws.Range("H92").Value = 410.33334
$ws.Range("I92").Value = 192.4
$ws.Range("K92").Value = 577.2
$ws.Range("M92").Value = 670.8

$ws.Range("H131").Value = 12659.421
$ws.Range("I131").Value = 773
$ws.Range("J131").Value = 21304.092
$ws.Range("K131").Value = 2319
$ws.Range("L131").Value = 63912.276
$ws.Range("M131").Value = 2721
$ws.Range("N131").Value = -73992.276

$ws.Range("H135").Value = 673.875
$ws.Range("I135").Value = 324.33334
$ws.Range("K135").Value = 2919.00006
$ws.Range("M135").Value = -384.0000600000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 919582.0600000001
$ws.Range("I11").Value = 1011400.3
$ws.Range("J11").Value = 1400
$ws.Range("K11").Value = 1011400.3
$ws.Range("L11").Value = 1400
$ws.Range("M11").Value = -1011261.3
$ws.Range("N11").Value = -1678

$ws.Range("H12").Value = 12500
$ws.Range("J12").Value = 12500
$ws.Range("L12").Value = 12500
$ws.Range("N12").Value = -12780

$ws.Range("H132").Value = 259561.4
$ws.Range("I132").Value = 315145.2
$ws.Range("J132").Value = 5464.143
$ws.Range("K132").Value = 945435.6000000001
$ws.Range("L132").Value = 16392.429
$ws.Range("M132").Value = -942905.6000000001
$ws.Range("N132").Value = -21452.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 250
$ws.Range("J14").Value = 250
$ws.Range("L14").Value = 250
$ws.Range("N14").Value = -594

$ws.Range("H40").Value = 3517.4375
$ws.Range("I40").Value = 3023.3333
$ws.Range("J40").Value = 4999.75
$ws.Range("K40").Value = 3023.3333
$ws.Range("L40").Value = 4999.75
$ws.Range("M40").Value = -2887.3333
$ws.Range("N40").Value = -5271.75

$ws.Range("H46").Value = 2369.1667
$ws.Range("I46").Value = 1738.6666
$ws.Range("J46").Value = 2999.6667
$ws.Range("K46").Value = 1738.6666
$ws.Range("L46").Value = 2999.6667
$ws.Range("M46").Value = -1550.6666
$ws.Range("N46").Value = -3375.6667

$ws.Range("H136").Value = 11824.25
$ws.Range("I136").Value = 4402
$ws.Range("J136").Value = 14298.333
$ws.Range("K136").Value = 13206
$ws.Range("L136").Value = 42894.999
$ws.Range("M136").Value = -10656
$ws.Range("N136").Value = -47994.999

$ws.Range("H137").Value = 59999
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 59999
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 59999
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -70199

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 5000
$ws.Range("K20").Value = 5000
$ws.Range("M20").Value = -4760

$ws.Range("H132").Value = 235442.72
$ws.Range("I132").Value = 273014.56
$ws.Range("K132").Value = 819043.6799999999
$ws.Range("M132").Value = -816513.6799999999

$ws.Range("H136").Value = 4185.7744
$ws.Range("I136").Value = 3173.2727
$ws.Range("J136").Value = 6660.778
$ws.Range("K136").Value = 9519.8181
$ws.Range("L136").Value = 19982.334
$ws.Range("M136").Value = -6969.8181
$ws.Range("N136").Value = -25082.334
